$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the comma after "ticket-blocking management"
# ------------------------------------------------------------------
$found1 = $d.Content.Find.Execute(
    "ticket-blocking management, and logs viewing pages",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ticket-blocking management and logs viewing pages",
    2)
if (-not $found1) {
    Write-Host "WARNING: step 1 text not found"
}

# ------------------------------------------------------------------
# 2) Rework the "force the reader to stop working ... no Internet
#    connectivity is available." sentence.
# ------------------------------------------------------------------
$found2 = $d.Content.Find.Execute(
    "force the reader to stop working. It will work as usual if the key expires and no Internet connectivity is available.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "stop the reader. It will work as usual if the key expires and there is no Internet connection.",
    2)
if (-not $found2) {
    Write-Host "WARNING: step 2 text not found"
}

# ------------------------------------------------------------------
# 3) Add a new bullet about the logs viewing page's log-clearing
#    feature, right after the "... Internet is OK." paragraph.
# ------------------------------------------------------------------
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd().EndsWith("Internet is OK.")) {
        $target = $para
    }
}

if ($target -eq $null) {
    Write-Host "WARNING: step 3 anchor paragraph not found"
} else {
    $target.Range.InsertParagraphAfter()
    $newPara = $target.Next()
    $newPara.Range.Text = "The logs viewing page supports clearing logs that are older than the specified minute (GDPR needs)."
}
